## TC11_Canine_Filter_Breed-Bouvier.xlsx - "Fixed ICDC breed all testcases"
##
## The CasesTab/SamplesTab/FilesTab rows on the "startup" sheet each carried
## a shared "StatQuery" (column C) that summarised counts for the selected
## breed. That Cypher query is replaced here with a corrected query that
## also reports Program/Study counts and splits file counts into
## "Case Files" vs "Study Files".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- 1. Replace the StatQuery (column C) used by all three tab rows -------
$newStatQuery = @"
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Bouvier des Flandres']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS ``Case Files``,
    count(distinct sf) AS ``Study Files``
"@
# the rest of the workbook's multi-line cells use CRLF line breaks
$newStatQuery = $newStatQuery -replace "`n", "`r`n"

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# --- 2. View: zoomed from 55% to 85%, whole sheet selected ----------------
$excel.ActiveWindow.Zoom = 85
$ws.Cells.Select()

# --- 3. Column widths settled to their post-edit autofit values -----------
$ws.Columns.Item(1).ColumnWidth = 11.333333333333334
$ws.Columns.Item(2).ColumnWidth = 72.33333333333333
$ws.Columns.Item(3).ColumnWidth = 61.0
$ws.Columns.Item(4).ColumnWidth = 46.333333333333336
$ws.Columns.Item(5).ColumnWidth = 45.333333333333336
